$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in header "occurenceID" -> "occurrenceID"
$ws.Range("C1").Value = "occurrenceID"

# Fix typo in header "scentificName" -> "scientificName"
$ws.Range("J1").Value = "scientificName"

# Move selection to J1 (last edited cell)
$ws.Range("J1").Select()
